# NSIT_YR_FIN.xlsx update:
# A new reporting period column is inserted before column D (the most recent
# fiscal-year column). All existing data in columns D:K shifts right to E:L,
# and the new column D is populated with the latest period's figures. A
# handful of prior-period cells in column E (rows 89, 101, 102) are also
# restated with slightly different values as part of this update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column at D; existing D:K shifts to E:L, preserving
#    the values/styles of the shifted cells.
$ws.Columns.Item(4).Insert()

# 2) The freshly inserted column D cells come in with "General" format.
#    Copy the (now-shifted) formatting from column E back onto column D so
#    the new column matches the date / number styling used throughout the
#    sheet.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Populate the new column D with the latest period's reported figures.
$newD = @{
    7   = 43465
    8   = 7080100
    9   = 6086400
    10  = 993700
    14  = 3700
    17  = 6846700
    18  = 233500
    20  = 1200
    21  = 272200
    22  = 22800
    23  = 211900
    24  = 53800
    26  = 158100
    27  = 158100
    29  = 5600
    32  = -1200
    33  = 163700
    35  = 163700
    38  = 43465
    41  = 142700
    43  = 1931700
    44  = 148500
    45  = 115700
    46  = 2338600
    48  = 73000
    49  = 279000
    52  = 85400
    54  = 2775900
    57  = 1282200
    58  = 1400
    59  = 253000
    60  = 1536700
    61  = 195500
    62  = 56800
    66  = 1789000
    72  = 704700
    76  = 987000
    80  = 43465
    81  = 163700
    83  = 37500
    89  = 292600
    91  = -17300
    94  = -91700
    100 = -159000
    101 = -5100
    102 = 36800
}
foreach ($r in $newD.Keys) {
    $ws.Cells.Item($r, 4).Value2 = $newD[$r]
}

# "Research Development" has no reported figure for any period, including
# the new one - carry the same "NA" marker used across the rest of the row.
$ws.Cells.Item(12, 4).Value2 = "NA"

# 4) A few of the prior-period figures (now in column E after the shift)
#    were restated along with this update.
$ws.Cells.Item(89, 5).Value2 = -307100
$ws.Cells.Item(101, 5).Value2 = 16100
$ws.Cells.Item(102, 5).Value2 = -98500
